$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.341.23'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '2.070.86'
$ws.Range("E3").Value = '  +4.97%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '''235.49'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("D6").Value = '''0.609'
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '''56.67'
$ws.Range("E8").Value = '  +5.42%  '
$ws.Range("D9").Value = '''0.379'
$ws.Range("E9").Value = '  +3.04%  '
$ws.Range("D10").Value = '''58.20'
$ws.Range("E10").Value = '  -2.64%  '
$ws.Range("D11").Value = '''0.0760'
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("E12").Value = '  +3.27%  '
$ws.Range("D13").Value = '2.365.11'
$ws.Range("E13").Value = '  +4.45%  '
$ws.Range("D14").Value = '''14.54'
$ws.Range("E14").Value = '  +5.23%  '
$ws.Range("D15").Value = '''21.00'
$ws.Range("E15").Value = '  +1.79%  '
$ws.Range("D16").Value = '''0.774'
$ws.Range("E16").Value = '  +3.79%  '
$ws.Range("D17").Value = '''5.24'
$ws.Range("E17").Value = '  +4.82%  '
$ws.Range("D18").Value = '2.043.24'
$ws.Range("E18").Value = '  +3.51%  '
$ws.Range("D19").Value = '37.254.36'
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").Value = '''5.89'
$ws.Range("E20").Value = '  +19.51%  '
$ws.Range("D21").Value = '''68.41'
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").Value = '0.0₃0809'
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").Value = '''223.77'
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("E25").Value = '  +3.25%  '
$ws.Range("D26").Value = '''2.43'
$ws.Range("E26").Value = '  +3.54%  '
$ws.Range("D27").Value = '''163.62'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").Value = '''8.86'
$ws.Range("E28").Value = '  +3.25%  '
$ws.Range("D29").Value = '''0.131'
$ws.Range("E29").Value = '  +6.52%  '
$ws.Range("D30").Value = '''19.30'
$ws.Range("E30").Value = '  +1.92%  '
$ws.Range("D31").Value = '''1.37'
$ws.Range("E31").Value = '  +7.31%  '
$ws.Range("E32").Value = '  +0.89%  '
$ws.Range("E33").Value = '  +1.51%  '
$ws.Range("D34").Value = '''0.0619'
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("D35").Value = '''2.53'
$ws.Range("E35").Value = '  +9.84%  '
$ws.Range("D36").Value = '''4.38'
$ws.Range("E36").Value = '  +4.07%  '
$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").Value = '''5.98'
$ws.Range("E37").Value = '  +16.07%  '
$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").Value = '''0.999'
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").Value = '''2.96'
$ws.Range("D42").Value = '''0.0949'
$ws.Range("E42").Value = '  +7.00%  '
$ws.Range("D43").Value = '1.477.51'
$ws.Range("E43").Value = '  +4.61%  '
$ws.Range("D44").Value = '''4.35'
$ws.Range("E44").Value = '  +17.03%  '
$ws.Range("D45").Value = '''94.82'
$ws.Range("E45").Value = '  +9.03%  '
$ws.Range("D46").Value = '''0.0208'
$ws.Range("E46").Value = '  +3.00%  '
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("D48").Value = '''16.01'
$ws.Range("E48").Value = '  +6.04%  '
$ws.Range("E49").Value = '  +3.47%  '
$ws.Range("E50").Value = '  +8.75%  '
$ws.Range("D51").Value = '''2.93'
$ws.Range("E51").Value = '  +2.70%  '
